$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registro_Errores")

# Widen column D (4) from 41 to 50
# Note: Excel's ColumnWidth property (character units) round-trips through an
# internal pixel-based representation that adds ~0.8333 when written back to
# the stored OOXML "width" attribute, so we compensate to land exactly on 50.
$ws.Columns.Item(4).ColumnWidth = 49.166666666666664

# New rows of data to append starting at row 8
$newRows = @(
    @(2, "2025-11-13 00:22:47", "sintoma_resuelto", "Virus ransomware eliminado de x_virus.exe", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_infectado_detectado", "Archivo spy_tool.exe puesto en cuarentena (Virus: spyware)", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_limpio_eliminado", "Error: kernel32.dll era un archivo limpio y fue eliminado", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_infectado_detectado", "Archivo adware_bundle.exe puesto en cuarentena (Virus: adware)", "No", "No", "No", 0),
    @(2, "2025-11-13 00:22:47", "archivo_limpio_cuarentena", "Falso positivo: logfile.log era seguro pero fue puesto en cuarentena", "No", "No", "No", 0)
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
